$wb = $excel.ActiveWorkbook

# Rename worksheets
$wsA = $wb.Worksheets.Item("double_barrel_type_a")
$wsA.Name = "double-barrel-A"

$wsB = $wb.Worksheets.Item("double_barrel_type_b")
$wsB.Name = "double-barrel-B"

# Update sheet A view: zoom, selection, not tab-selected
$wsA.Activate()
$wsA.Range("C16").Select()
$excel.ActiveWindow.Zoom = 290

# Update sheet B view: zoom, selection, and becomes the active tab
$wsB.Activate()
$wsB.Range("B2").Select()
$excel.ActiveWindow.Zoom = 216
